$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire rows for the two records that were dropped from the table
# ("RM 232" at row 26 and "SC 92" at row 28). Delete the lower row first so
# the row index for the upper one stays valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# After the deletions the remaining rows shift up to close the gaps. Apply
# the individual cell value updates (re-rolled missing/present values) using
# the new row numbers.
$ws.Range("F2").Value = $null

$ws.Range("F5").Value = 17.66

$ws.Range("D6").Value = -14.2
$ws.Range("F6").Value = 16.43

$ws.Range("D8").Value = $null

$ws.Range("F10").Value = $null

$ws.Range("D12").Value = -14.1

$ws.Range("F13").Value = $null

$ws.Range("D14").Value = $null

$ws.Range("D17").Value = -14.7

$ws.Range("D18").Value = -15.2

$ws.Range("D19").Value = $null

$ws.Range("D20").Value = $null

$ws.Range("D23").Value = -13.9

$ws.Range("F24").Value = 16.78

$ws.Range("C27").Value = 10
$ws.Range("D27").Value = $null

$ws.Range("F28").Value = $null

$ws.Range("C29").Value = $null

$ws.Range("F30").Value = 16.89

$ws.Range("C32").Value = $null
